$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 8: add the date of the event and its description
# Copy the date formatting used by the row above onto A8 before setting the value
$ws.Range("A7").Copy()
$ws.Range("A8").PasteSpecial([Microsoft.Office.Interop.Excel.XlPasteType]::xlPasteFormats)
$excel.CutCopyMode = $false

$ws.Range("A8").Value = 43929
$ws.Range("B8").Value = "Rendu de la version de la version 1.0 du projet et fin du sprint 6"

# Update the selection to B8
$ws.Range("B8").Select()
